$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text: "MODEL_CONDITION" -> "MODELCONDITION"
$ws.Range("E1").Value = "MODELCONDITION"

# Drop the redundant leading column (its values duplicate the last column)
# and shift everything else (B:F) one column to the left (A:E).
$ws.Columns.Item(1).Delete()
